# Revert "Drop in all data files from 3.0 RMI script"
# Set the Boolean Exempt Process Emissions From Carbon Tax control lever
# back to 1 (on the "BEPEfCT" worksheet, cell B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 1
